$d = $word.ActiveDocument

# ----------------------------------------------------------------------
# Edit 1: "...then on "Proteins" (on the right side of newly opened page)
# and then on to "InterPro" (it's the sixth entry under Data Resources)."
# becomes:
# "...then type "InterPro" into the box "find a data resource or tool"
# and click "Search". From the returned results, click "InterPro" to
# switch to Interpro Website."
# ----------------------------------------------------------------------
$d.Content.Find.Execute(", then on “Proteins” (on the right side of newly opened page) and then on to “", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    ", then type “", 2) | Out-Null

$d.Content.Find.Execute("” (it’s the sixth entry under Data Resources).", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "” into the box “find a data resource or tool” and click “Search”. From the returned results, click “InterPro” to switch to Interpro Website. ", 2) | Out-Null

# ----------------------------------------------------------------------
# Edit 2: comment on "How many transmembrane sections do you find now?"
# ----------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("How many transmembrane sections do you find now?", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cmt1 = $d.Comments.Add($rng1, "10 ? ")
$cmt1.Author = "Tao Fang (tfang)"
$cmt1.Initial = "TF("

# ----------------------------------------------------------------------
# Edit 3: "PDBe" -> "PDB"
# ----------------------------------------------------------------------
$d.Content.Find.Execute("PDBe", $true, $false, $false, $false, $false, $true, 1, $false, "PDB", 2) | Out-Null

# ----------------------------------------------------------------------
# Edit 4: comment on ", even viruses"
# ----------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(", even viruses", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cmt2 = $d.Comments.Add($rng2, "Now I didn’t see the virus option ")
$cmt2.Author = "Tao Fang (tfang)"
$cmt2.Initial = "TF("

# ----------------------------------------------------------------------
# Edit 5: comment on "of" in "towards the middle of the page"
# ----------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("towards the middle of the page", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sub3 = $d.Range($rng3.Start, $rng3.End)
$sub3.Find.Execute("of", $false, $true, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$cmt3Text = "Now website is “ pfam-legacy.xfam.org ” its said  remain available until January 2023 , but  I can stil use it until 09052023 at least "
$cmt3 = $d.Comments.Add($sub3, $cmt3Text)
$cmt3.Author = "Tao Fang (tfang)"
$cmt3.Initial = "TF("

Write-Output "all edits applied"
